# Generate Report for Handback
# Updates the zh-cn and de-de localization-status sheets to reflect that the
# handback has completed: status text changes, "Latest Target File" /
# "Latest Handback File" columns (E/F) get populated with hyperlinks, and the
# "Latest Handback DateTime" column (G) gets the real handback timestamp.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

function Update-LangSheet {
    param($sheetName, $xlfName, $handbackDatetime, $handoffUrlPrefix, $handbackUrlPrefix)

    $ws = $wb.Worksheets.Item($sheetName)

    # --- Row 2 (a.md) ---------------------------------------------------
    $ws.Range("B2").Value = $statusText

    $ws.Hyperlinks.Add($ws.Range("E2"), ($handoffUrlPrefix + "/e2e/a.md"), "", "", "a.md")

    $ws.Hyperlinks.Add($ws.Range("F2"), ($handbackUrlPrefix + "/" + $xlfName), "", "", $xlfName)

    $ws.Range("G2").Value = $handbackDatetime

    # --- Row 3 (b.md) ----------------------------------------------------
    $ws.Range("B3").Value = $statusText

    # NB: mirrors the existing (pre-existing) data quirk in this sheet where
    # row 3's handoff-file hyperlink (C3) already pointed at the "a." xlf
    # instead of a "b." one -- the new Target/Handback columns follow suit.
    $ws.Hyperlinks.Add($ws.Range("E3"), ($handoffUrlPrefix + "/e2e/a.md"), "", "", "a.md")

    $ws.Hyperlinks.Add($ws.Range("F3"), ($handbackUrlPrefix + "/" + $xlfName), "", "", $xlfName)

    $ws.Range("G3").Value = $handbackDatetime
}

$handoffUrlPrefix = "https://github.com/OpenLocalizationTest/oltest/blob/1384363847824c1cd1144a21cee68ba9f7fb5ecf"

Update-LangSheet "zh-cn" `
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" `
    "2016-02-18 07:46:35" `
    $handoffUrlPrefix `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/400b3126d75fdb0e856c0766fc27404cd071642b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/hb"

Update-LangSheet "de-de" `
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" `
    "2016-02-18 07:46:55" `
    $handoffUrlPrefix `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dfdf1b3f892b8b425a83bb8ef215555391766ab0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/hb"

Write-Host "Handback report generated."
